$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-obsolete rows 12 and 13 first (MuSCs sending-cluster group dropped),
# so remaining rows 2-11 can be updated in place afterwards.
$ws.Rows("12:13").Delete()

# Update data rows 2-11 with recalculated TPM values and corrected sending/target cluster assignments
# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Hc"
$ws.Range("C2").Value = "C5ar2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1183866666666667
$ws.Range("H2").Value = 0.35516
$ws.Range("I2").Value = 0.7705291475929099
$ws.Range("J2").Value = 0.77052914759291
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.120074
$ws.Range("N2").Value = 0.360222
$ws.Range("O2").Value = 0.00855715934293867
$ws.Range("P2").Value = 0.008580002469517843
$ws.Range("Q2").Value = 0.01421516061333333
$ws.Range("R2").Value = 0.12793644552
$ws.Range("S2").Value = 0.006593540694331239
$ws.Range("T2").Value = 0.006611141989182646

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Hc"
$ws.Range("C3").Value = "C5ar2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1183866666666667
$ws.Range("H3").Value = 0.35516
$ws.Range("I3").Value = 0.7705291475929099
$ws.Range("J3").Value = 0.77052914759291
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.4465006666666667
$ws.Range("N3").Value = 1.339502
$ws.Range("O3").Value = 0.03182018881185777
$ws.Range("P3").Value = 0.03190513202392993
$ws.Range("Q3").Value = 0.0528597255911111
$ws.Range("R3").Value = 0.4757375303199999
$ws.Range("S3").Value = 0.02451838296144622
$ws.Range("T3").Value = 0.02458383418223798

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Hc"
$ws.Range("C4").Value = "C5ar2"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1183866666666667
$ws.Range("H4").Value = 0.35516
$ws.Range("I4").Value = 0.7705291475929099
$ws.Range("J4").Value = 0.77052914759291
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.033632000000001
$ws.Range("N4").Value = 21.100896
$ws.Range("O4").Value = 0.5012568064992622
$ws.Range("P4").Value = 0.5025948992261415
$ws.Range("Q4").Value = 0.83268824704
$ws.Range("R4").Value = 7.49419422336
$ws.Range("S4").Value = 0.3862329798370207
$ws.Range("T4").Value = 0.3872640192852633

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Hc"
$ws.Range("C5").Value = "C5ar2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1183866666666667
$ws.Range("H5").Value = 0.35516
$ws.Range("I5").Value = 0.7705291475929099
$ws.Range("J5").Value = 0.77052914759291
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.112075
$ws.Range("N5").Value = 0.22415
$ws.Range("O5").Value = 0.007987104896645831
$ws.Range("P5").Value = 0.005338950851259568
$ws.Range("Q5").Value = 0.01326818566666667
$ws.Range("R5").Value = 0.07960911399999999
$ws.Range("S5").Value = 0.006154297127747669
$ws.Range("T5").Value = 0.004113817248461477

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Hc"
$ws.Range("C6").Value = "C5ar2"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1183866666666667
$ws.Range("H6").Value = 0.35516
$ws.Range("I6").Value = 0.7705291475929099
$ws.Range("J6").Value = 0.77052914759291
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.319711333333333
$ws.Range("N6").Value = 18.959134
$ws.Range("O6").Value = 0.4503787404492957
$ws.Range("P6").Value = 0.4515810154291511
$ws.Range("Q6").Value = 0.7481695590488888
$ws.Range("R6").Value = 6.733526031439999
$ws.Range("S6").Value = 0.3470299469723642
$ws.Range("T6").Value = 0.3479563348877645

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Hc"
$ws.Range("C7").Value = "C5ar2"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.03525666666666667
$ws.Range("H7").Value = 0.10577
$ws.Range("I7").Value = 0.22947085240709
$ws.Range("J7").Value = 0.22947085240709
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.120074
$ws.Range("N7").Value = 0.360222
$ws.Range("O7").Value = 0.00855715934293867
$ws.Range("P7").Value = 0.008580002469517843
$ws.Range("Q7").Value = 0.004233408993333333
$ws.Range("R7").Value = 0.03810068094000001
$ws.Range("S7").Value = 0.001963618648607431
$ws.Range("T7").Value = 0.001968860480335197

# Row 8
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Hc"
$ws.Range("C8").Value = "C5ar2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.03525666666666667
$ws.Range("H8").Value = 0.10577
$ws.Range("I8").Value = 0.22947085240709
$ws.Range("J8").Value = 0.22947085240709
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4465006666666667
$ws.Range("N8").Value = 1.339502
$ws.Range("O8").Value = 0.03182018881185777
$ws.Range("P8").Value = 0.03190513202392993
$ws.Range("Q8").Value = 0.01574212517111111
$ws.Range("R8").Value = 0.14167912654
$ws.Range("S8").Value = 0.007301805850411552
$ws.Range("T8").Value = 0.007321297841691947

# Row 9
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Hc"
$ws.Range("C9").Value = "C5ar2"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.03525666666666667
$ws.Range("H9").Value = 0.10577
$ws.Range("I9").Value = 0.22947085240709
$ws.Range("J9").Value = 0.22947085240709
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.033632000000001
$ws.Range("N9").Value = 21.100896
$ws.Range("O9").Value = 0.5012568064992622
$ws.Range("P9").Value = 0.5025948992261415
$ws.Range("Q9").Value = 0.24798241888
$ws.Range("R9").Value = 2.23184176992
$ws.Range("S9").Value = 0.1150238266622415
$ws.Range("T9").Value = 0.1153308799408782

# Row 10
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Hc"
$ws.Range("C10").Value = "C5ar2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.03525666666666667
$ws.Range("H10").Value = 0.10577
$ws.Range("I10").Value = 0.22947085240709
$ws.Range("J10").Value = 0.22947085240709
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.112075
$ws.Range("N10").Value = 0.22415
$ws.Range("O10").Value = 0.007987104896645831
$ws.Range("P10").Value = 0.005338950851259568
$ws.Range("Q10").Value = 0.003951390916666667
$ws.Range("R10").Value = 0.0237083455
$ws.Range("S10").Value = 0.001832807768898161
$ws.Range("T10").Value = 0.001225133602798092

# Row 11
$ws.Range("A11").Value = "Inflammatory-Mac"
$ws.Range("B11").Value = "Hc"
$ws.Range("C11").Value = "C5ar2"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.03525666666666667
$ws.Range("H11").Value = 0.10577
$ws.Range("I11").Value = 0.22947085240709
$ws.Range("J11").Value = 0.22947085240709
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 6.319711333333333
$ws.Range("N11").Value = 18.959134
$ws.Range("O11").Value = 0.4503787404492957
$ws.Range("P11").Value = 0.4515810154291511
$ws.Range("Q11").Value = 0.2228119559088889
$ws.Range("R11").Value = 2.00530760318
$ws.Range("S11").Value = 0.1033487934769314
$ws.Range("T11").Value = 0.1036246805413866
